$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1780343333333333
$ws.Range("H2").Value = 0.534103
$ws.Range("I2").Value = 0.003649670474736916
$ws.Range("J2").Value = 0.003649670474736915
$ws.Range("O2").Value = 0.0002880357555630755
$ws.Range("P2").Value = 0.0002880357555630755
$ws.Range("Q2").Value = 0.002170238523333333
$ws.Range("R2").Value = 0.01953214671
$ws.Range("S2").Value = 0.000001051235592747096
$ws.Range("T2").Value = 0.000001051235592747096

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1780343333333333
$ws.Range("H3").Value = 0.534103
$ws.Range("I3").Value = 0.003649670474736916
$ws.Range("J3").Value = 0.003649670474736915
$ws.Range("M3").Value = 32.87103466666667
$ws.Range("N3").Value = 98.61310400000001
$ws.Range("O3").Value = 0.7767049471988007
$ws.Range("P3").Value = 0.7767049471988008
$ws.Range("Q3").Value = 5.852172742856888
$ws.Range("R3").Value = 52.669554685712
$ws.Range("S3").Value = 0.002834717113373558
$ws.Range("T3").Value = 0.002834717113373558

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1780343333333333
$ws.Range("H4").Value = 0.534103
$ws.Range("I4").Value = 0.003649670474736916
$ws.Range("J4").Value = 0.003649670474736915
$ws.Range("M4").Value = 9.43791
$ws.Range("N4").Value = 28.31373
$ws.Range("O4").Value = 0.2230070170456362
$ws.Range("P4").Value = 0.2230070170456362
$ws.Range("Q4").Value = 1.68027201491
$ws.Range("R4").Value = 15.12244813419
$ws.Range("S4").Value = 0.0008139021257706108
$ws.Range("T4").Value = 0.0008139021257706106

# Row 5
$ws.Range("I5").Value = 0.09908483984804967
$ws.Range("J5").Value = 0.09908483984804965
$ws.Range("O5").Value = 0.0002880357555630755
$ws.Range("P5").Value = 0.0002880357555630755
$ws.Range("S5").Value = 0.00002853997671047932
$ws.Range("T5").Value = 0.00002853997671047931

# Row 6
$ws.Range("I6").Value = 0.09908483984804967
$ws.Range("J6").Value = 0.09908483984804965
$ws.Range("M6").Value = 32.87103466666667
$ws.Range("N6").Value = 98.61310400000001
$ws.Range("O6").Value = 0.7767049471988007
$ws.Range("P6").Value = 0.7767049471988008
$ws.Range("Q6").Value = 158.8805353806347
$ws.Range("R6").Value = 1429.924818425712
$ws.Range("S6").Value = 0.07695968530238104
$ws.Range("T6").Value = 0.07695968530238104

# Row 7
$ws.Range("I7").Value = 0.09908483984804967
$ws.Range("J7").Value = 0.09908483984804965
$ws.Range("M7").Value = 9.43791
$ws.Range("N7").Value = 28.31373
$ws.Range("O7").Value = 0.2230070170456362
$ws.Range("P7").Value = 0.2230070170456362
$ws.Range("Q7").Value = 45.61767552741
$ws.Range("R7").Value = 410.55907974669
$ws.Range("S7").Value = 0.02209661456895815
$ws.Range("T7").Value = 0.02209661456895815

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 35.63223
$ws.Range("H8").Value = 106.89669
$ws.Range("I8").Value = 0.7304540385283456
$ws.Range("J8").Value = 0.7304540385283456
$ws.Range("O8").Value = 0.0002880357555630755
$ws.Range("P8").Value = 0.0002880357555630755
$ws.Range("Q8").Value = 0.4343568837
$ws.Range("R8").Value = 3.9092119533
$ws.Range("S8").Value = 0.0002103968808916119
$ws.Range("T8").Value = 0.0002103968808916119

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 35.63223
$ws.Range("H9").Value = 106.89669
$ws.Range("I9").Value = 0.7304540385283456
$ws.Range("J9").Value = 0.7304540385283456
$ws.Range("M9").Value = 32.87103466666667
$ws.Range("N9").Value = 98.61310400000001
$ws.Range("O9").Value = 0.7767049471988007
$ws.Range("P9").Value = 0.7767049471988008
$ws.Range("Q9").Value = 1171.26826758064
$ws.Range("R9").Value = 10541.41440822576
$ws.Range("S9").Value = 0.5673472654263093
$ws.Range("T9").Value = 0.5673472654263094

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 35.63223
$ws.Range("H10").Value = 106.89669
$ws.Range("I10").Value = 0.7304540385283456
$ws.Range("J10").Value = 0.7304540385283456
$ws.Range("M10").Value = 9.43791
$ws.Range("N10").Value = 28.31373
$ws.Range("O10").Value = 0.2230070170456362
$ws.Range("P10").Value = 0.2230070170456362
$ws.Range("Q10").Value = 336.2937798393
$ws.Range("R10").Value = 3026.6440185537
$ws.Range("S10").Value = 0.1628963762211446
$ws.Range("T10").Value = 0.1628963762211446

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5521946666666667
$ws.Range("H11").Value = 1.656584
$ws.Range("I11").Value = 0.011319887201011
$ws.Range("J11").Value = 0.011319887201011
$ws.Range("O11").Value = 0.0002880357555630755
$ws.Range("P11").Value = 0.0002880357555630755
$ws.Range("Q11").Value = 0.006731252986666667
$ws.Range("R11").Value = 0.06058127688
$ws.Range("S11").Value = 0.000003260532262831992
$ws.Range("T11").Value = 0.000003260532262831992

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5521946666666667
$ws.Range("H12").Value = 1.656584
$ws.Range("I12").Value = 0.011319887201011
$ws.Range("J12").Value = 0.011319887201011
$ws.Range("M12").Value = 32.87103466666667
$ws.Range("N12").Value = 98.61310400000001
$ws.Range("O12").Value = 0.7767049471988007
$ws.Range("P12").Value = 0.7767049471988008
$ws.Range("Q12").Value = 18.15121003074844
$ws.Range("R12").Value = 163.360890276736
$ws.Range("S12").Value = 0.008792212390757632
$ws.Range("T12").Value = 0.008792212390757632

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5521946666666667
$ws.Range("H13").Value = 1.656584
$ws.Range("I13").Value = 0.011319887201011
$ws.Range("J13").Value = 0.011319887201011
$ws.Range("M13").Value = 9.43791
$ws.Range("N13").Value = 28.31373
$ws.Range("O13").Value = 0.2230070170456362
$ws.Range("P13").Value = 0.2230070170456362
$ws.Range("Q13").Value = 5.211563566480001
$ws.Range("R13").Value = 46.90407209832
$ws.Range("S13").Value = 0.00252441427799054
$ws.Range("T13").Value = 0.00252441427799054

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 7.585023666666667
$ws.Range("H14").Value = 22.755071
$ws.Range("I14").Value = 0.155491563947857
$ws.Range("J14").Value = 0.1554915639478569
$ws.Range("O14").Value = 0.0002880357555630755
$ws.Range("P14").Value = 0.0002880357555630755
$ws.Range("Q14").Value = 0.09246143849666667
$ws.Range("R14").Value = 0.83215294647
$ws.Range("S14").Value = 0.00004478713010540525
$ws.Range("T14").Value = 0.00004478713010540524

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 7.585023666666667
$ws.Range("H15").Value = 22.755071
$ws.Range("I15").Value = 0.155491563947857
$ws.Range("J15").Value = 0.1554915639478569
$ws.Range("M15").Value = 32.87103466666667
$ws.Range("N15").Value = 98.61310400000001
$ws.Range("O15").Value = 0.7767049471988007
$ws.Range("P15").Value = 0.7767049471988008
$ws.Range("Q15").Value = 249.3275758944871
$ws.Range("R15").Value = 2243.948183050384
$ws.Range("S15").Value = 0.1207710669659792
$ws.Range("T15").Value = 0.1207710669659792

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 7.585023666666667
$ws.Range("H16").Value = 22.755071
$ws.Range("I16").Value = 0.155491563947857
$ws.Range("J16").Value = 0.1554915639478569
$ws.Range("M16").Value = 9.43791
$ws.Range("N16").Value = 28.31373
$ws.Range("O16").Value = 0.2230070170456362
$ws.Range("P16").Value = 0.2230070170456362
$ws.Range("Q16").Value = 71.58677071387001
$ws.Range("R16").Value = 644.28093642483
$ws.Range("S16").Value = 0.03467570985177237
$ws.Range("T16").Value = 0.03467570985177237
